# precios4.xlsx - simplify the long product-description labels in column A
# and move the active selection, matching the author's "Add files via upload" edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Shorten the six shared product descriptions -------------------------
# Rows 2-40  : "INTEGRA 6000 ... SIN ANTENA"  -> "INTEGRA 6000 SIN ANTENA"
$ws.Range("A2:A40").Value = "INTEGRA 6000 SIN ANTENA"

# Rows 41-79 : "INTEGRA 6000 ... CON ANTENA"  -> "INTEGRA 6000 CON ANTENA"
$ws.Range("A41:A79").Value = "INTEGRA 6000 CON ANTENA"

# Rows 80-118: "PANTALLA i6000 ... SIN ANTENA" -> "PANTALLA i6000 SIN ANTENA"
$ws.Range("A80:A118").Value = "PANTALLA i6000 SIN ANTENA"

# Rows 119-157: "PANTALLA i6000 ... CON ANTENA" -> "PANTALLA i6000 CON ANTENA"
$ws.Range("A119:A157").Value = "PANTALLA i6000 CON ANTENA"

# Rows 158-196: "Equipo de Corte por Surco con cualquier modelo..." -> "TEKMATIC"
$ws.Range("A158:A196").Value = "TEKMATIC"

# Rows 197-235: "Equipo de Corte por Surco con GATILLO..." -> "GATILLO PARA DOSIFICADOR MATERMACC"
$ws.Range("A197:A235").Value = "GATILLO PARA DOSIFICADOR MATERMACC"

# --- Move the sheet's active selection (as captured when the file was saved) ---
$ws.Range("A200").Select()
